$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Finish C4-Algorithm row (row 9) ---
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = "YES"
$ws.Range("E9").Value = "DONE"

# --- Insert new row for C5-ImproveRSSI ---
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "C5-ImproveRSSI"
$ws.Range("D10").Value = "NO"
$ws.Range("E10").Value = "NOTHING"

# --- Fix up the SUM row (now row 11) ---
$ws.Range("C11").Formula = "=SUM(C2:C10)"

# Selection per the diff
$ws.Range("C12").Select()
